$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric / reused-string cells first (do not affect sharedStrings order) ---
$ws.Range("A7").Value = 1.6
$ws.Range("A8").Value = 1.7
$ws.Range("A9").Value = 1.8
$ws.Range("F7").Value = "ок"
$ws.Range("F8").Value = "ок"
$ws.Range("F9").Value = "ок"

# --- New text cells, set in the exact order the original author typed them ---
# (B7, D7, B8, D8, E7, E8, B9, D9, E9) so sharedStrings are appended in matching order
$ws.Range("B7").Value = "Попытка создания проекта при незаполненом поле `"Сокращенное название`""
$ws.Range("D7").Value = "1. Заполнить поле `"Название`", например значением `"я люблю сокращать названия`"`n2. Нажать кнопку сохранить"
$ws.Range("B8").Value = "Попытка создания проекта при заполнении поля `"Название`" невалидным значением"
$ws.Range("D8").Value = "1. Заполнить поле `"Название`" невалидным значением, например, `"Проеееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееект 1`"`n2. Заполнить поле `"Сокращенное название`" допустимым значением, например `"Сокращение 1`"`n3. Нажать кнопку сохранить"
$ws.Range("E7").Value = "отображение алерта об незаполненном обязательно поле `"Сокращенное название`""
$ws.Range("E8").Value = "отображение алерта об неверном заполнении обязательного поля  `"Название`""
$ws.Range("B9").Value = "Попытка создания проекта при заполнении поля `"Сокращенное название`" невалидным значением"
$ws.Range("D9").Value = "1. Заполнить поле `"Сокращенное название`" невалидным значением, например, `"Проеееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееееект 1`"`n2. Заполнить поле `"Сокращенное название`" допустимым значением, например `"Сокращение 1`"`n3. Нажать кнопку сохранить"
$ws.Range("E9").Value = "отображение алерта об неверном заполнении обязательного поля  `"Сокращенное название`""

# --- Formatting: wrap text + top alignment for new data cells (matches existing columns) ---
$ws.Range("A7:E9").WrapText = $true
$ws.Range("A7:E9").VerticalAlignment = -4160

# --- Row heights ---
$ws.Rows(7).RowHeight = 77
$ws.Rows(8).RowHeight = 255
$ws.Rows(9).RowHeight = 255

# --- Extend merged precondition cell down through the new rows ---
$ws.Range("C2:C9").Merge()

# --- Apply the "Good" (green) cell style to the whole Result column, as already used for "ок" ---
$ws.Range("F2:F9").Style = "Good"

# --- Column D width tweak ---
$ws.Columns("D").ColumnWidth = 39.3

# --- Sheet view: zoom + selection ---
$ws.Application.ActiveWindow.Zoom = 89
$ws.Range("A2:B9").Select()

Write-Output "done"
